# feat: Context Menu First Implementation
#
# Inserts a new "UI View" worksheet immediately before the existing
# "Tiles View" worksheet, and populates it with the new UI prefab rows
# (ContextMenuContainer / ButtonPrefab) using the same 3-column layout
# ("IDS" | "Architecture ID" | "Prefab resource path") as the other
# "* View" sheets in this workbook.

$wb = $excel.ActiveWorkbook

# Anchor on the existing "Tiles View" sheet so the new sheet is inserted
# right before it (tab order: ... , "Tile Types", "UI View", "Tiles View").
$tilesView = $wb.Worksheets.Item("Tiles View")

$uiView = $wb.Worksheets.Add($tilesView)
$uiView.Name = "UI View"

# Header row - same headers used by "Prefabs View" / "Tiles View".
$uiView.Range("A1").Value = "IDS"
$uiView.Range("B1").Value = "Architecture ID"
$uiView.Range("C1").Value = "Prefab resource path"

# ContextMenuContainer prefab entry.
$uiView.Range("A2").Value = "MENU_CONTAINER"
$uiView.Range("B2").Value = "ContextMenuContainer"
$uiView.Range("C2").Value = "Prefabs/UI/ContextMenuContainer"

# ButtonPrefab entry.
$uiView.Range("A3").Value = "MENU_BUTTON"
$uiView.Range("B3").Value = "ButtonPrefab"
$uiView.Range("C3").Value = "Prefabs/UI/ButtonPrefab"

# Keep the originally active sheet selected (Worksheets.Add activates the
# newly created sheet by default).
$wb.Worksheets.Item(1).Activate()
